$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2021年" row (row 5) below the existing 2018-2020 rows.
# Copy the header-cell formatting from A4 (centered, bold, bordered) onto A5,
# then set its label.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "2021年"

# Copy AG4's formatting/empty-string placeholder onto AG5 (AG column is
# blank for every year in this sheet).
$ws.Range("AG4").Copy()
$ws.Range("AG5").PasteSpecial(-4122)
$ws.Range("AG5").Value = ""

# Fill in the 2021 growth figures for every other indicator column (B:DK,
# excluding AG).
$values = @{
    "B5" = 2.8
    "C5" = 23.9
    "D5" = -28.9
    "E5" = 37.7
    "F5" = 5.2
    "G5" = 81.2
    "H5" = 13
    "I5" = 5.9
    "J5" = 6.3
    "K5" = -7.9
    "L5" = 2.3
    "M5" = 22.5
    "N5" = -9
    "O5" = -0.8
    "P5" = 5.9
    "Q5" = -8.5
    "R5" = 29.8
    "S5" = 12.4
    "T5" = 12
    "U5" = 10.9
    "V5" = 5
    "W5" = 18.9
    "X5" = 14.4
    "Y5" = 15.1
    "Z5" = 32.2
    "AA5" = 9.4
    "AB5" = 19.9
    "AC5" = 25.1
    "AD5" = 7.4
    "AE5" = 14.2
    "AF5" = -6.8
    "AH5" = -22.2
    "AI5" = -4.6
    "AJ5" = -18.6
    "AK5" = 16.1
    "AL5" = 3.8
    "AM5" = 1.4
    "AN5" = -9.6
    "AO5" = -10.7
    "AP5" = -27.2
    "AQ5" = 6.8
    "AR5" = 4.7
    "AS5" = 26
    "AT5" = 36.5
    "AU5" = 49.7
    "AV5" = 4.7
    "AW5" = 18.2
    "AX5" = -6.1
    "AY5" = -5.8
    "AZ5" = 10.6
    "BA5" = 1.2
    "BB5" = 5.5
    "BC5" = 11.8
    "BD5" = 66.9
    "BE5" = 6.5
    "BF5" = -1
    "BG5" = 14.5
    "BH5" = -3.1
    "BI5" = -3.4
    "BJ5" = 13.4
    "BK5" = 19.4
    "BL5" = -0.7
    "BM5" = 2.7
    "BN5" = -3.8
    "BO5" = -2.4
    "BP5" = 5.1
    "BQ5" = -13.9
    "BR5" = 10.5
    "BS5" = -2.8
    "BT5" = -2.9
    "BU5" = 14.1
    "BV5" = 4.7
    "BW5" = 10.5
    "BX5" = 23.5
    "BY5" = 20.7
    "BZ5" = 2.5
    "CA5" = 18.3
    "CB5" = 85.9
    "CC5" = 15
    "CD5" = -26.1
    "CE5" = -3.1
    "CF5" = 12.1
    "CG5" = 14.5
    "CH5" = -22.7
    "CI5" = 12.4
    "CJ5" = -6.8
    "CK5" = 11.8
    "CL5" = 4.4
    "CM5" = -18.7
    "CN5" = 16.5
    "CO5" = 10.8
    "CP5" = 22.7
    "CQ5" = 13.5
    "CR5" = -27.2
    "CS5" = 15.1
    "CT5" = 9.9
    "CU5" = 13.2
    "CV5" = 2.6
    "CW5" = 19.2
    "CX5" = 17
    "CY5" = 21.6
    "CZ5" = 43.4
    "DA5" = 11.6
    "DB5" = 3.9
    "DC5" = 25.9
    "DD5" = 17.4
    "DE5" = -5.5
    "DF5" = 14.2
    "DG5" = 25.3
    "DH5" = 10.5
    "DI5" = 9.3
    "DJ5" = 11.1
    "DK5" = 28.3
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
